# Update Supplemental Table 2 with the finalized title/caption text and
# formatting, per "Added the updated version of Supplemental Table 2."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: table title (merged A1:D1) -------------------------------
$ws.Range("A1:D1").Value = "Table S2: Distribution of unique TIR, LRR, and IG domain isoforms across sequenced coral genomes"
$ws.Range("A1:D1").WrapText = $true
$ws.Rows(1).RowHeight = 27

# --- Row 2: table caption / methods note (merged A2:D2) --------------
$ws.Range("A2:D2").Value = "Unique isofom copies of TIR, LRR, and Ig domains as annotated using Transdecoder to predict open reading frames followed by HMMER searched for the domains in the prediced reading frames."
$ws.Range("A2:D2").VerticalAlignment = -4108
$ws.Rows(2).RowHeight = 48.6

# --- Selection, to match the saved view state -------------------------
$ws.Range("A2:D2").Select()
